# Extend the CDISC CT starter-kit workbook with the VS/AE/EX/DS/DM codelists
# (VS Test Code, VS Test Name, No/Yes Response, Severity, Route of
# Administration, Disposition Category, Epoch, Dose Frequency) and their terms.
#
# Sheet "Codelists"       (ws1) -> one header row per codelist (rows 7-14)
# Sheet "Codelists_terms" (ws2) -> one row per codelist term   (rows 18-47)
#
# Cells that are blank in the target (no TERM_CODE / DECODE) are still
# "touched" via .Style = 'Normal' so the cell is materialized (present, but
# empty) rather than omitted from the row, matching the source data model.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Codelists): append codelist header rows 7-14 ---
$ws1.Range('A7').Value = 'Y'
$ws1.Range('B7').Value = 'VS Test Code'
$ws1.Range('C7').Value = 'C66741'
$ws1.Range('D7').Value = 'VSTESTCD'
$ws1.Range('E7').Value = 'Yes'

$ws1.Range('A8').Value = 'Y'
$ws1.Range('B8').Value = 'VS Test Name'
$ws1.Range('C8').Value = 'C67153'
$ws1.Range('D8').Value = 'VSTEST'
$ws1.Range('E8').Value = 'Yes'

$ws1.Range('A9').Value = 'Y'
$ws1.Range('B9').Value = 'No/Yes Response'
$ws1.Range('C9').Value = 'C66742'
$ws1.Range('D9').Value = 'NY'
$ws1.Range('E9').Value = 'No'

$ws1.Range('A10').Value = 'Y'
$ws1.Range('B10').Value = 'Severity'
$ws1.Range('C10').Value = 'C66769'
$ws1.Range('D10').Value = 'AESEV'
$ws1.Range('E10').Value = 'No'

$ws1.Range('A11').Value = 'Y'
$ws1.Range('B11').Value = 'Route of Administration'
$ws1.Range('C11').Value = 'C66729'
$ws1.Range('D11').Value = 'ROUTE'
$ws1.Range('E11').Value = 'Yes'

$ws1.Range('A12').Value = 'Y'
$ws1.Range('B12').Value = 'Disposition Category'
$ws1.Range('C12').Value = 'C66767'
$ws1.Range('D12').Value = 'DSCAT'
$ws1.Range('E12').Value = 'Yes'

$ws1.Range('A13').Value = 'Y'
$ws1.Range('B13').Value = 'Epoch'
$ws1.Range('C13').Value = 'C71738'
$ws1.Range('D13').Value = 'EPOCH'
$ws1.Range('E13').Value = 'Yes'

$ws1.Range('A14').Value = 'Y'
$ws1.Range('B14').Value = 'Dose Frequency'
$ws1.Range('C14').Value = 'C71113'
$ws1.Range('D14').Value = 'FREQ'
$ws1.Range('E14').Value = 'Yes'

# --- Sheet2 (Codelists_terms): append codelist term rows 18-47 ---
$ws2.Range('A18').Value = 'Y'
$ws2.Range('B18').Value = 'VS Test Code'
$ws2.Range('C18').Value = 'VSTESTCD'
$ws2.Cells.Item(18,4).Style = 'Normal'
$ws2.Range('E18').Value = 'SYSBP'
$ws2.Cells.Item(18,6).Style = 'Normal'

$ws2.Range('A19').Value = 'Y'
$ws2.Range('B19').Value = 'VS Test Code'
$ws2.Range('C19').Value = 'VSTESTCD'
$ws2.Cells.Item(19,4).Style = 'Normal'
$ws2.Range('E19').Value = 'DIABP'
$ws2.Cells.Item(19,6).Style = 'Normal'

$ws2.Range('A20').Value = 'Y'
$ws2.Range('B20').Value = 'VS Test Code'
$ws2.Range('C20').Value = 'VSTESTCD'
$ws2.Cells.Item(20,4).Style = 'Normal'
$ws2.Range('E20').Value = 'HR'
$ws2.Cells.Item(20,6).Style = 'Normal'

$ws2.Range('A21').Value = 'Y'
$ws2.Range('B21').Value = 'VS Test Code'
$ws2.Range('C21').Value = 'VSTESTCD'
$ws2.Cells.Item(21,4).Style = 'Normal'
$ws2.Range('E21').Value = 'WEIGHT'
$ws2.Cells.Item(21,6).Style = 'Normal'

$ws2.Range('A22').Value = 'Y'
$ws2.Range('B22').Value = 'VS Test Code'
$ws2.Range('C22').Value = 'VSTESTCD'
$ws2.Cells.Item(22,4).Style = 'Normal'
$ws2.Range('E22').Value = 'HEIGHT'
$ws2.Cells.Item(22,6).Style = 'Normal'

$ws2.Range('A23').Value = 'Y'
$ws2.Range('B23').Value = 'VS Test Name'
$ws2.Range('C23').Value = 'VSTEST'
$ws2.Cells.Item(23,4).Style = 'Normal'
$ws2.Range('E23').Value = 'Systolic Blood Pressure'
$ws2.Cells.Item(23,6).Style = 'Normal'

$ws2.Range('A24').Value = 'Y'
$ws2.Range('B24').Value = 'VS Test Name'
$ws2.Range('C24').Value = 'VSTEST'
$ws2.Cells.Item(24,4).Style = 'Normal'
$ws2.Range('E24').Value = 'Diastolic Blood Pressure'
$ws2.Cells.Item(24,6).Style = 'Normal'

$ws2.Range('A25').Value = 'Y'
$ws2.Range('B25').Value = 'VS Test Name'
$ws2.Range('C25').Value = 'VSTEST'
$ws2.Cells.Item(25,4).Style = 'Normal'
$ws2.Range('E25').Value = 'Heart Rate'
$ws2.Cells.Item(25,6).Style = 'Normal'

$ws2.Range('A26').Value = 'Y'
$ws2.Range('B26').Value = 'VS Test Name'
$ws2.Range('C26').Value = 'VSTEST'
$ws2.Cells.Item(26,4).Style = 'Normal'
$ws2.Range('E26').Value = 'Weight'
$ws2.Cells.Item(26,6).Style = 'Normal'

$ws2.Range('A27').Value = 'Y'
$ws2.Range('B27').Value = 'VS Test Name'
$ws2.Range('C27').Value = 'VSTEST'
$ws2.Cells.Item(27,4).Style = 'Normal'
$ws2.Range('E27').Value = 'Height'
$ws2.Cells.Item(27,6).Style = 'Normal'

$ws2.Range('A28').Value = 'Y'
$ws2.Range('B28').Value = 'No/Yes Response'
$ws2.Range('C28').Value = 'NY'
$ws2.Range('D28').Value = 'C49488'
$ws2.Range('E28').Value = 'N'
$ws2.Range('F28').Value = 'No'

$ws2.Range('A29').Value = 'Y'
$ws2.Range('B29').Value = 'No/Yes Response'
$ws2.Range('C29').Value = 'NY'
$ws2.Range('D29').Value = 'C49487'
$ws2.Range('E29').Value = 'Y'
$ws2.Range('F29').Value = 'Yes'

$ws2.Range('A30').Value = 'Y'
$ws2.Range('B30').Value = 'Severity'
$ws2.Range('C30').Value = 'AESEV'
$ws2.Cells.Item(30,4).Style = 'Normal'
$ws2.Range('E30').Value = 'MILD'
$ws2.Range('F30').Value = 'Mild'

$ws2.Range('A31').Value = 'Y'
$ws2.Range('B31').Value = 'Severity'
$ws2.Range('C31').Value = 'AESEV'
$ws2.Cells.Item(31,4).Style = 'Normal'
$ws2.Range('E31').Value = 'MODERATE'
$ws2.Range('F31').Value = 'Moderate'

$ws2.Range('A32').Value = 'Y'
$ws2.Range('B32').Value = 'Severity'
$ws2.Range('C32').Value = 'AESEV'
$ws2.Cells.Item(32,4).Style = 'Normal'
$ws2.Range('E32').Value = 'SEVERE'
$ws2.Range('F32').Value = 'Severe'

$ws2.Range('A33').Value = 'Y'
$ws2.Range('B33').Value = 'Route of Administration'
$ws2.Range('C33').Value = 'ROUTE'
$ws2.Cells.Item(33,4).Style = 'Normal'
$ws2.Range('E33').Value = 'ORAL'
$ws2.Range('F33').Value = 'Oral'

$ws2.Range('A34').Value = 'Y'
$ws2.Range('B34').Value = 'Route of Administration'
$ws2.Range('C34').Value = 'ROUTE'
$ws2.Cells.Item(34,4).Style = 'Normal'
$ws2.Range('E34').Value = 'INTRAVENOUS'
$ws2.Range('F34').Value = 'Intravenous'

$ws2.Range('A35').Value = 'Y'
$ws2.Range('B35').Value = 'Route of Administration'
$ws2.Range('C35').Value = 'ROUTE'
$ws2.Cells.Item(35,4).Style = 'Normal'
$ws2.Range('E35').Value = 'SUBCUTANEOUS'
$ws2.Range('F35').Value = 'Subcutaneous'

$ws2.Range('A36').Value = 'Y'
$ws2.Range('B36').Value = 'Disposition Category'
$ws2.Range('C36').Value = 'DSCAT'
$ws2.Cells.Item(36,4).Style = 'Normal'
$ws2.Range('E36').Value = 'PROTOCOL MILESTONE'
$ws2.Range('F36').Value = 'Protocol Milestone'

$ws2.Range('A37').Value = 'Y'
$ws2.Range('B37').Value = 'Disposition Category'
$ws2.Range('C37').Value = 'DSCAT'
$ws2.Cells.Item(37,4).Style = 'Normal'
$ws2.Range('E37').Value = 'DISPOSITION EVENT'
$ws2.Range('F37').Value = 'Disposition Event'

$ws2.Range('A38').Value = 'Y'
$ws2.Range('B38').Value = 'Epoch'
$ws2.Range('C38').Value = 'EPOCH'
$ws2.Cells.Item(38,4).Style = 'Normal'
$ws2.Range('E38').Value = 'SCREENING'
$ws2.Range('F38').Value = 'Screening'

$ws2.Range('A39').Value = 'Y'
$ws2.Range('B39').Value = 'Epoch'
$ws2.Range('C39').Value = 'EPOCH'
$ws2.Cells.Item(39,4).Style = 'Normal'
$ws2.Range('E39').Value = 'TREATMENT'
$ws2.Range('F39').Value = 'Treatment'

$ws2.Range('A40').Value = 'Y'
$ws2.Range('B40').Value = 'Epoch'
$ws2.Range('C40').Value = 'EPOCH'
$ws2.Cells.Item(40,4).Style = 'Normal'
$ws2.Range('E40').Value = 'FOLLOW-UP'
$ws2.Range('F40').Value = 'Follow-up'

$ws2.Range('A41').Value = 'Y'
$ws2.Range('B41').Value = 'Dose Frequency'
$ws2.Range('C41').Value = 'FREQ'
$ws2.Cells.Item(41,4).Style = 'Normal'
$ws2.Range('E41').Value = 'QD'
$ws2.Range('F41').Value = 'Every Day'

$ws2.Range('A42').Value = 'Y'
$ws2.Range('B42').Value = 'Dose Frequency'
$ws2.Range('C42').Value = 'FREQ'
$ws2.Cells.Item(42,4).Style = 'Normal'
$ws2.Range('E42').Value = 'BID'
$ws2.Range('F42').Value = 'Twice a Day'

$ws2.Range('A43').Value = 'Y'
$ws2.Range('B43').Value = 'Dose Frequency'
$ws2.Range('C43').Value = 'FREQ'
$ws2.Cells.Item(43,4).Style = 'Normal'
$ws2.Range('E43').Value = 'TID'
$ws2.Range('F43').Value = 'Three Times a Day'

$ws2.Range('A44').Value = 'Y'
$ws2.Range('B44').Value = 'Dose Frequency'
$ws2.Range('C44').Value = 'FREQ'
$ws2.Cells.Item(44,4).Style = 'Normal'
$ws2.Range('E44').Value = 'QID'
$ws2.Range('F44').Value = 'Four Times a Day'

$ws2.Range('A45').Value = 'Y'
$ws2.Range('B45').Value = 'Dose Frequency'
$ws2.Range('C45').Value = 'FREQ'
$ws2.Cells.Item(45,4).Style = 'Normal'
$ws2.Range('E45').Value = 'Q2H'
$ws2.Range('F45').Value = 'Every 2 Hours'

$ws2.Range('A46').Value = 'Y'
$ws2.Range('B46').Value = 'Dose Frequency'
$ws2.Range('C46').Value = 'FREQ'
$ws2.Cells.Item(46,4).Style = 'Normal'
$ws2.Range('E46').Value = 'Q4H'
$ws2.Range('F46').Value = 'Every 4 Hours'

$ws2.Range('A47').Value = 'Y'
$ws2.Range('B47').Value = 'Dose Frequency'
$ws2.Range('C47').Value = 'FREQ'
$ws2.Cells.Item(47,4).Style = 'Normal'
$ws2.Range('E47').Value = 'Q6H'
$ws2.Range('F47').Value = 'Every 6 Hours'
